$d = $word.ActiveDocument

# 1) Insert new opening sentence "Es ist 00:20 Uhr. " before the Foyer paragraph.
$d.Content.Find.Execute(
    "Das Foyer des Hotel Aurora ist erfüllt von aufgeregten Stimmen.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Es ist 00:20 Uhr. Das Foyer des Hotel Aurora ist erfüllt von aufgeregten Stimmen.",
    2)

# 2) Rework the "Doch draußen..." paragraph: add the "denn durch die Wolken..." clause
#    and change "hat" -> "hatte" / "Sturm" -> "Schnee".
$d.Content.Find.Execute(
    "Doch draußen, im peitschenden Schneetreiben, war kaum noch etwas zu erkennen. Die Nacht war so finster wie die Stimmung der Gäste. Der Wind hat die meisten Spuren verweht, der Sturm hat jede klare Kontur verschluckt. Nur wenige Hinweise lassen sich überhaupt noch deuten:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Doch draußen, im peitschenden Schneetreiben, war kaum noch etwas zu erkennen, denn durch die Wolken war die Nacht war so finster wie die Stimmung der Gäste. Der Wind hatte die meisten Spuren verweht, der Schnee jede klare Kontur verschluckt. Nur wenige Hinweise lassen sich überhaupt noch deuten:",
    2)

# 3) Drop the "Und wer zum Teufel ... Schuhe unterwegs?" aside (and the stray space before "Sie").
$d.Content.Find.Execute(
    ", oder hatte der Wind sie im Schnee verzerrt? Und wer zum Teufel ist bei dem Wetter hier draußen ohne Schuhe unterwegs? Sie",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", oder hatte der Wind sie im Schnee verzerrt? Sie",
    2)

# 4) Remove the ellipsis and add "gab es einen Gehilfen" before the question mark.
$d.Content.Find.Execute(
    "Einer allein könnte das kaum schaffen. Oder…?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Einer allein könnte das kaum schaffen. Oder gab es einen Gehilfen?",
    2)

# 5) "das" -> "dass" (also removes the now-unneeded grammar-check markers around it).
$d.Content.Find.Execute(
    "ihre Stimme nun fest, das selbst sie nicht verbergen kann:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ihre Stimme nun fest, dass selbst sie nicht verbergen kann:",
    2)

# 6) Merge the three short closing paragraphs into one and change "haben begonnen" -> "können beginnen".
$d.Content.Find.Execute(
    "Das Hotel Aurora wirkt plötzlich kleiner.^pViel kleiner.^pDie Ermittlungen haben begonnen.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Das Hotel Aurora wirkt plötzlich kleiner. Viel kleiner. Die Ermittlungen können beginnen.",
    2)
